$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1 - copy style from G1 (header row style) and set value
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Save column values mirror the inverse-ish pattern described in the diff
$saveValues = @(1, 0, 1, 1, 1, 1, 0, 0, 0, 1, 1, 0, 1)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
